# The post originally at row 219 ("「漢字の成り立ち」...") was removed from the
# workbook. Deleting the entire row shifts every subsequent row up by one,
# which matches the target diff (old row 220 becomes the new row 219, etc.)
# and shrinks the used range from A1:C394 to A1:C393.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(219).Delete()
